# The post that used to live at row 346 ("「私たちは異なる。だから調和する」")
# was removed from the source data; delete its entire row here so every
# subsequent row shifts up by one and the sheet's used range shrinks from
# A1:C489 to A1:C488.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(346).Delete()
